$d = $word.ActiveDocument
$gitPara = $d.Paragraphs.Item(9)
$gitRange = $gitPara.Range
$insertionPoint = $d.Range($gitRange.Start, $gitRange.Start)
$insertionPoint.InsertBefore("React`r" + "flowbite-react.com`r" + "`r")

# Re-fetch by fresh index lookups (not stale object refs)
$bulletTemplate = $d.Paragraphs.Item(13)  # "git init" now at index 13
Write-Output ("template text: " + $bulletTemplate.Range.Text)
$listTemplate = $bulletTemplate.Range.ListFormat.ListTemplate

$flowbitePara = $d.Paragraphs.Item(10)
Write-Output ("flowbite text before: [" + $flowbitePara.Range.Text + "]")
$flowbitePara.Style = "List Paragraph"
$flowbitePara.Range.ListFormat.ApplyListTemplateWithLevel($listTemplate, $true, 1, $false, 1)

$blankPara = $d.Paragraphs.Item(11)
Write-Output ("blank text before: [" + $blankPara.Range.Text + "]")
$blankPara.Style = "List Paragraph"
$blankPara.Range.ListFormat.ApplyListTemplateWithLevel($listTemplate, $true, 1, $false, 1)

for ($i=1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    Write-Output ($i.ToString() + ": [" + $p.Range.Text + "] style=" + $p.Style.NameLocal + " listType=" + $p.Range.ListFormat.ListType)
}
